$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.939.41"
$ws.Range("E2").Value = "  +4.96%  "

$ws.Range("D3").Value = "2.268.56"
$ws.Range("E3").Value = "  +2.19%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.92%  "

$ws.Range("E7").Value = "  +3.35%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  +4.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0799"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.59%  "

$ws.Range("E13").Value = "  +2.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.60%  "

$ws.Range("D15").Value = "2.617.54"
$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("E16").Value = "  +2.70%  "

$ws.Range("D17").Value = "2.254.92"
$ws.Range("E17").Value = "  +3.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.761"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.93%  "

$ws.Range("D19").Value = "41.860.78"
$ws.Range("E19").Value = "  +4.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.14%  "

$ws.Range("E21").Value = "  +2.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.07"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "242.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.33%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  +3.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.18%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.42%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.41%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("E34").Value = "  +4.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0747"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.07"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.57%  "

$ws.Range("E39").Value = "  +4.90%  "

$ws.Range("E40").Value = "  +3.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.11%  "

$ws.Range("E42").Value = "  +6.13%  "

$ws.Range("D43").Value = "2.065.86"
$ws.Range("E43").Value = "  -1.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.86%  "

$ws.Range("E45").Value = "  +3.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.90%  "

$ws.Range("E47").Value = "  +7.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.85%  "

$ws.Range("E49").Value = "  +4.00%  "

$ws.Range("E50").Value = "  +3.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.67%  "
